# softwarestack.pptx — "Update the figure of softwarestack"
#
# 1) The auto-updating "datetimeFigureOut" footer field (cached as
#    2017/4/7 on the slide master and all 11 slide layouts) is refreshed
#    to 2017/6/25.
# 2) On slide 1, the "I/O Cluster #0-1" shape's label is split so the
#    "#0-1" suffix becomes its own run (text unchanged).
# 3) On slide 1, the "Compute Cluster #0-15 " shape's label is split so
#    the suffix becomes its own run, updated from "#0-15 " to "#1-16 ".

$p = $ppt.ActivePresentation

# --- 1) Refresh the cached "update automatically" date field -------------
function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2017/4/7") {
                $tr.Text = "2017/6/25"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DateField $master.CustomLayouts.Item($li).Shapes
}

# --- 2) & 3) Split the cluster-count suffix into its own run -------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if (-not $sh.HasTextFrame) {
        continue
    }

    $tr = $sh.TextFrame.TextRange
    for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
        $para = $tr.Paragraphs($pi)
        $text = $para.Text

        if ($text -eq "I/O Cluster #0-1") {
            $suffix = $para.Characters(13, 4)
            $suffix.Text = "#0-1"
        }
        elseif ($text -eq "Compute Cluster #0-15 ") {
            $suffix = $para.Characters(17, 6)
            $suffix.Text = "#1-16 "
        }
    }
}
